# Apply updated crypto price/volume figures to Sheet1 (columns D and E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.623.56"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "'2.544.99"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'311.97"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "'99.91"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "'35.95"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'2.939.39"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "'16.20"
$ws.Range("E15").Value = "  +7.30%  "
$ws.Range("D16").Value = "'2.616.58"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "'0.840"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "'42.606.32"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "'12.30"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "'0.0₃0950"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "'69.03"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'242.51"
$ws.Range("E23").Value = "  -4.35%  "
$ws.Range("D24").Value = "'2.90"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").Value = "'26.35"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("D29").Value = "'40.00"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").Value = "'158.55"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  +14.20%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("D36").Value = "'2.04"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").Value = "'3.16"
$ws.Range("E37").Value = "  -5.35%  "
$ws.Range("D38").Value = "'17.95"
$ws.Range("E38").Value = "  -7.03%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "'4.23"
$ws.Range("E41").Value = "  +10.43%  "
$ws.Range("D42").Value = "'21.70"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").Value = "'3.34"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "'1.959.15"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "'8.92"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'2.796.62"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "'80.83"
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "'72.32"
$ws.Range("E51").Value = "  -2.85%  "
